$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D:E data range so numeric-looking strings
# (e.g. "313.65") are preserved as text, matching the inlineStr cell type
# used by the source workbook instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.074.14'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.874.90'
$ws.Range('E3').Value = '  -1.62%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').Value = '313.65'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('D7').Value = '0.5081'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '0.3849'
$ws.Range('E8').Value = '  -2.14%  '
$ws.Range('D9').Value = '0.09013'
$ws.Range('E9').Value = '  -2.97%  '
$ws.Range('D10').Value = '1.123'
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('D11').Value = '41.58'
$ws.Range('E11').Value = '  -0.74%  '
$ws.Range('D12').Value = '6.343'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('D14').Value = '1.865.51'
$ws.Range('E14').Value = '  -2.25%  '
$ws.Range('D15').Value = '7.197'
$ws.Range('E15').Value = '  -1.41%  '
$ws.Range('E16').Value = '  +0.37%  '
$ws.Range('D17').Value = '0.00001111'
$ws.Range('E17').Value = '  -0.89%  '
$ws.Range('D19').Value = '0.06596'
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('D20').Value = '18.17'
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('D22').Value = '6.112'
$ws.Range('E22').Value = '  -1.81%  '
$ws.Range('D23').Value = '28.090.90'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').Value = '11.42'
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').Value = '2.281'
$ws.Range('E25').Value = '  -1.50%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '2.538'
$ws.Range('E26').Value = '  -3.62%  '
$ws.Range('D27').Value = '2.084.02'
$ws.Range('E27').Value = '  -1.90%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '157.75'
$ws.Range('E28').Value = '  +0.24%  '
$ws.Range('D29').Value = '20.78'
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '126.55'
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.1052'
$ws.Range('E31').Value = '  -1.63%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '1.060'
$ws.Range('E32').Value = '  -2.83%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '5.617'
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '3.604'
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').Value = '9.655'
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.06578'
$ws.Range('E36').Value = '  -1.58%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.02428'
$ws.Range('E37').Value = '  +0.87%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = '0.2178'
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '1.208'
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '1.264'
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.6394'
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '11.45'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('B43').Value = 'InternetComputer(DFINITY)'
$ws.Range('C43').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D43').Value = '4.910'
$ws.Range('E43').Value = '  -1.75%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').Value = '0.6016'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '13.15'
$ws.Range('E45').Value = '  -0.94%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.672'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('B47').Value = 'WEMIXTOKEN'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '1.276'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D48').Value = '1.240'
$ws.Range('E48').Value = '  +5.06%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.994'
$ws.Range('E49').Value = '  -1.17%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '121.45'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '79.69'
$ws.Range('E51').Value = '  +1.61%  '

# Restore default (unstyled) cell style so no stray number-format style
# index is left attached to the cells, matching the original file layout.
$ws.Range("D2:E51").Style = "Normal"
